# Generate Report for Handoff
#
# Updates the localization-status report to reflect that zh-cn / de-de
# have moved from "In Translation" to "Ready for handoff", and refreshes
# the associated handoff timestamps. Also widens the "Status" columns
# (Overview!E:F and the per-locale sheets' column C) to fit the new,
# longer status text - mirroring the column auto-resize that Excel would
# perform when the report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "In Translation"
$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------------
# Row 2: File Name | Path And Name | Extension | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-15 14:55:54"

# --- zh-cn sheet ------------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-15 14:55:50"

# --- de-de sheet ------------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-15 14:55:54"

# --- Widen the Status columns to fit "Ready for handoff" -------------------
# (13.4101845877511 -> 17.2159881591797 chars in the authored report)
$newStatusColWidth = 16.333333333333336
$wsOverview.Range("E1:F1").ColumnWidth = $newStatusColWidth
$wsZhCn.Range("C1").ColumnWidth = $newStatusColWidth
$wsDeDe.Range("C1").ColumnWidth = $newStatusColWidth
